$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-03-14 Thursday" "2024-03-15 Friday"

Replace-Text "817×8=" "781×3="
Replace-Text "795×9=" "911×2="
Replace-Text "243×3=" "959×6="
Replace-Text "832×8=" "963×4="
Replace-Text "531×5=" "206×2="

Replace-Text "705×4=" "166×6="
Replace-Text "705×9=" "230×4="
Replace-Text "162×5=" "756×6="
Replace-Text "727×2=" "966×2="
Replace-Text "983×5=" "341×2="

Replace-Text "716×9=" "201×7="
Replace-Text "935×8=" "803×7="
Replace-Text "969×8=" "543×9="
Replace-Text "699×5=" "275×2="
Replace-Text "291×5=" "287×8="

Replace-Text "232×7=" "693×8="
Replace-Text "236×4=" "585×4="
Replace-Text "157×3=" "932×6="
Replace-Text "419×9=" "154×8="
Replace-Text "746×5=" "522×2="

Replace-Text "808×5=" "900×5="
Replace-Text "200×4=" "478×6="
Replace-Text "863×4=" "301×8="
Replace-Text "662×2=" "382×6="
Replace-Text "562×2=" "364×4="
